$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new row of data (the 1s/2s Cl entries) was inserted above the existing
# table, pushing the header (row 3) and all data rows (4-28) down by one
# (header -> row 4, data -> rows 5-29).
$ws.Rows(3).Insert()

# Update the active selection to match the post-edit view state.
$ws.Range("E7").Select()
